# Remove the trailing "Ver no Jupiter ..." / copyright footer block that
# Jekyll appended right after the "Requisitos" section, along with the
# blank paragraph that separated it from the requirement text.
#
# Resulting structure: the "LOQ4044: ..." paragraph is followed directly
# by the blank paragraph that used to sit right before the page-break
# paragraph at the end of the document.

$d = $word.ActiveDocument

# Locate the first paragraph of the footer block ("Ver no Jupiter ...").
$r1 = $d.Content
$found1 = $r1.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", $false, $true, $false, $false, $false, $true, 1, $false, "", 0)

if ($found1) {
    $verParaIndex = $r1.Paragraphs.Item(1).Index

    # The empty paragraph right before it is also part of the block being
    # removed (it only existed to separate the requirement text from the
    # footer).
    $deleteStart = $d.Paragraphs.Item($verParaIndex - 1).Range.Start

    # Locate the last paragraph of the footer block (the copyright line).
    $r2 = $d.Content
    $found2 = $r2.Find.Execute("Powered by Jekyll and Github pages", $false, $true, $false, $false, $false, $true, 1, $false, "", 0)
    $copyParaIndex = $r2.Paragraphs.Item(1).Index
    $deleteEnd = $d.Paragraphs.Item($copyParaIndex).Range.End

    $deleteRange = $d.Range($deleteStart, $deleteEnd)
    $deleteRange.Delete()
}
